$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp in the title row (A1)
$ws.Range("A1").Value = "Datos actualizados a 24 de Mayo de 2020 a las 02:05"

# Row 4: Estados Unidos - refreshed totals
$ws.Range("B4").Value = 1666801
$ws.Range("C4").Value = 21902
$ws.Range("D4").Value = 446874
$ws.Range("E4").Value = 1121253
$ws.Range("G4").Value = 1027
$ws.Range("H4").Value = 98674

# Rows 120/121: Andorra and Uruguay swap places (with refreshed stats)
$ws.Range("A120").Value = "Uruguay"
$ws.Range("B120").Value = 764
$ws.Range("C120").Value = 11
$ws.Range("D120").Value = 616
$ws.Range("E120").Value = 126
$ws.Range("G120").Value = 2
$ws.Range("H120").Value = 22

$ws.Range("A121").Value = "Principado de Andorra"
$ws.Range("B121").Value = 762
$ws.Range("D121").Value = 653
$ws.Range("E121").Value = 58
$ws.Range("H121").Value = 51

# Row 166: Guyana - refreshed active/recovered counts
$ws.Range("D166").Value = 58
$ws.Range("E166").Value = 59

# Rows 170/171: Monaco and Bahamas swap places (with refreshed stats)
$ws.Range("A170").Value = "Bahamas"
$ws.Range("B170").Value = 100
$ws.Range("C170").Value = 3
$ws.Range("D170").Value = 45
$ws.Range("E170").Value = 44
$ws.Range("H170").Value = 11

$ws.Range("A171").Value = "Monaco"
$ws.Range("B171").Value = 98
$ws.Range("C171").Value = 1
$ws.Range("D171").Value = 90
$ws.Range("E171").Value = 4
$ws.Range("H171").Value = 4
